$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.05080186434227
$ws.Range("D2").Value = 1.050249810047854
$ws.Range("E2").Value = 1.060003990030804
$ws.Range("F2").Value = 1.067980466610889
$ws.Range("I2").Value = 1.046606682602516
$ws.Range("J2").Value = 1.05583309628923
$ws.Range("K2").Value = 1.053004610620832
$ws.Range("L2").Value = 1.06273195172191
$ws.Range("M2").Value = 1.070686875065925
$ws.Range("N2").Value = 1.05733250039171
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.052014061284733
$ws.Range("D3").Value = 1.051167977104799
$ws.Range("E3").Value = 1.061413608125844
$ws.Range("F3").Value = 1.069493681379185
$ws.Range("I3").Value = 1.047005694322364
$ws.Range("J3").Value = 1.056693839362022
$ws.Range("K3").Value = 1.053734913709601
$ws.Range("L3").Value = 1.06395440644704
$ws.Range("M3").Value = 1.072014246483254
$ws.Range("N3").Value = 1.058194465818394
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.05279746527675
$ws.Range("D4").Value = 1.051761189961945
$ws.Range("E4").Value = 1.062325463957048
$ws.Range("F4").Value = 1.070472704132239
$ws.Range("I4").Value = 1.047262032593623
$ws.Range("J4").Value = 1.057249304227319
$ws.Range("K4").Value = 1.054205938955303
$ws.Range("L4").Value = 1.064744627865078
$ws.Range("M4").Value = 1.072872506913524
$ws.Range("N4").Value = 1.058750719507513
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.053126579680875
$ws.Range("D5").Value = 1.052010362303743
$ws.Range("E5").Value = 1.062708750291446
$ws.Range("F5").Value = 1.070884259536921
$ws.Range("I5").Value = 1.047369355716905
$ws.Range("J5").Value = 1.057482466349722
$ws.Range("K5").Value = 1.054403593544523
$ws.Range("L5").Value = 1.065076652413702
$ws.Range("M5").Value = 1.073233171366599
$ws.Range("N5").Value = 1.058984212746873
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.053181826123727
$ws.Range("D6").Value = 1.052052186917682
$ws.Range("E6").Value = 1.062773102497405
$ws.Range("F6").Value = 1.070953360134082
$ws.Range("I6").Value = 1.04738734988637
$ws.Range("J6").Value = 1.057521594525255
$ws.Range("K6").Value = 1.054436759301476
$ws.Range("L6").Value = 1.065132390098241
$ws.Range("M6").Value = 1.073293719951136
$ws.Range("N6").Value = 1.059023396488904
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.052801863813898
$ws.Range("D7").Value = 1.051764520255037
$ws.Range("E7").Value = 1.062330585669284
$ws.Range("F7").Value = 1.070478203450725
$ws.Range("I7").Value = 1.047263468383749
$ws.Range("J7").Value = 1.057252421144365
$ws.Range("K7").Value = 1.054208581455329
$ws.Range("L7").Value = 1.064749065113169
$ws.Range("M7").Value = 1.072877326706551
$ws.Range("N7").Value = 1.058753840850938
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.051211734450589
$ws.Range("D8").Value = 1.05056029685345
$ws.Range("E8").Value = 1.060480434012661
$ws.Range("F8").Value = 1.06849189420825
$ws.Range("I8").Value = 1.046741914032007
$ws.Range("J8").Value = 1.056124299017303
$ws.Range("K8").Value = 1.053251737779101
$ws.Range("L8").Value = 1.063145251845415
$ws.Range("M8").Value = 1.071135601140844
$ws.Range("N8").Value = 1.057624116661058
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.04840216086682
$ws.Range("D9").Value = 1.048431309271082
$ws.Range("E9").Value = 1.057218018020967
$ws.Range("F9").Value = 1.064990557250941
$ws.Range("I9").Value = 1.045808654842678
$ws.Range("J9").Value = 1.054124865025325
$ws.Range("K9").Value = 1.051553860747625
$ws.Range("L9").Value = 1.060312898741372
$ws.Range("M9").Value = 1.068061376512633
$ws.Range("N9").Value = 1.055621843243473
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.046523821947692
$ws.Range("D10").Value = 1.047007169813101
$ws.Range("E10").Value = 1.055041290496678
$ws.Range("F10").Value = 1.062655208376574
$ws.Range("I10").Value = 1.045176850233391
$ws.Range("J10").Value = 1.052784008136719
$ws.Range("K10").Value = 1.050413896604615
$ws.Range("L10").Value = 1.058420208364607
$ws.Range("M10").Value = 1.066008198360382
$ws.Range("N10").Value = 1.054279082184287
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.04570917932395
$ws.Range("D11").Value = 1.046389333722281
$ws.Range("E11").Value = 1.054098254166484
$ws.Range("F11").Value = 1.061643639751161
$ws.Range("I11").Value = 1.044900970075839
$ws.Range("J11").Value = 1.052201495776108
$ws.Range("K11").Value = 1.049918344820675
$ws.Range("L11").Value = 1.057599536142916
$ws.Range("M11").Value = 1.065118210722145
$ws.Range("N11").Value = 1.053695742589308
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.045406383943557
$ws.Range("D12").Value = 1.046159663112978
$ws.Range("E12").Value = 1.053747887752961
$ws.Range("F12").Value = 1.061267839859075
$ws.Range("I12").Value = 1.044798148101218
$ws.Range("J12").Value = 1.051984834518383
$ws.Range("K12").Value = 1.049733981032742
$ws.Range("L12").Value = 1.057294527739424
$ws.Range("M12").Value = 1.064787481818916
$ws.Range("N12").Value = 1.053478773647746
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.045471343735009
$ws.Range("D13").Value = 1.046208936366233
$ws.Range("E13").Value = 1.053823046284209
$ws.Range("F13").Value = 1.061348452946586
$ws.Range("I13").Value = 1.044820219527706
$ws.Range("J13").Value = 1.052031322273991
$ws.Range("K13").Value = 1.04977354102382
$ws.Range("L13").Value = 1.057359961044173
$ws.Range("M13").Value = 1.064858431070424
$ws.Range("N13").Value = 1.0535253274213
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.045684154284207
$ws.Range("D14").Value = 1.046370352738705
$ws.Range("E14").Value = 1.054069294449813
$ws.Range("F14").Value = 1.06161257723886
$ws.Range("I14").Value = 1.044892477883438
$ws.Range("J14").Value = 1.052183592426429
$ws.Range("K14").Value = 1.049903111259407
$ws.Range("L14").Value = 1.057574327636982
$ws.Range("M14").Value = 1.065090875615893
$ws.Range("N14").Value = 1.053677813814819
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.04581524715105
$ws.Range("D15").Value = 1.046469782933028
$ws.Range("E15").Value = 1.054221005201078
$ws.Range("F15").Value = 1.061775305028711
$ws.Range("I15").Value = 1.044936952507389
$ws.Range("J15").Value = 1.052277372547375
$ws.Range("K15").Value = 1.049982904756884
$ws.Range("L15").Value = 1.057706382730357
$ws.Range("M15").Value = 1.065234072614674
$ws.Range("N15").Value = 1.053771727114294
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.046577859156755
$ws.Range("D16").Value = 1.047048148618233
$ws.Range("E16").Value = 1.055103865707847
$ws.Range("F16").Value = 1.062722335019831
$ws.Range("I16").Value = 1.045195110802697
$ws.Range("J16").Value = 1.05282262702261
$ws.Range("K16").Value = 1.050446743674323
$ws.Range("L16").Value = 1.058474649534571
$ws.Range("M16").Value = 1.066067243472078
$ws.Range("N16").Value = 1.054317755913425
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.047055872045591
$ws.Range("D17").Value = 1.047410626083674
$ws.Range("E17").Value = 1.055657523129514
$ws.Range("F17").Value = 1.063316284665921
$ws.Range("I17").Value = 1.045356428576699
$ws.Range("J17").Value = 1.053164136681127
$ws.Range("K17").Value = 1.050737176575496
$ws.Range("L17").Value = 1.058956258189832
$ws.Range("M17").Value = 1.066589611906102
$ws.Range("N17").Value = 1.054659750554829
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.047334562810269
$ws.Range("D18").Value = 1.047621939913688
$ws.Range("E18").Value = 1.05598041438696
$ws.Range("F18").Value = 1.063662692365394
$ws.Range("I18").Value = 1.045450300251818
$ws.Range("J18").Value = 1.053363149158417
$ws.Range("K18").Value = 1.050906394172738
$ws.Range("L18").Value = 1.059237064197321
$ws.Range("M18").Value = 1.066894209482621
$ws.Range("N18").Value = 1.054859045652664
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.047429567899397
$ws.Range("D19").Value = 1.047693973348274
$ws.Range("E19").Value = 1.056090504036216
$ws.Range("F19").Value = 1.063780803011976
$ws.Range("I19").Value = 1.04548227041604
$ws.Range("J19").Value = 1.053430976092051
$ws.Range("K19").Value = 1.050964061352745
$ws.Range("L19").Value = 1.059332793598215
$ws.Range("M19").Value = 1.066998054097836
$ws.Range("N19").Value = 1.054926968908324
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.047004598857127
$ws.Range("D20").Value = 1.047371747402967
$ws.Range("E20").Value = 1.055598125910207
$ws.Range("F20").Value = 1.063252562957491
$ws.Range("I20").Value = 1.045339143706826
$ws.Range("J20").Value = 1.053127515014785
$ws.Range("K20").Value = 1.050706035234107
$ws.Range("L20").Value = 1.058904597367001
$ws.Range("M20").Value = 1.066533576179926
$ws.Range("N20").Value = 1.05462307688152
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.045621492458819
$ws.Range("D21").Value = 1.046322824594267
$ws.Range("E21").Value = 1.053996782777673
$ws.Range("F21").Value = 1.061534800875494
$ws.Range("I21").Value = 1.044871209212257
$ws.Range("J21").Value = 1.052138760696958
$ws.Range("K21").Value = 1.049864964183582
$ws.Range("L21").Value = 1.057511206869078
$ws.Range("M21").Value = 1.065022430616535
$ws.Range("N21").Value = 1.05363291841915
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.044750714534912
$ws.Range("D22").Value = 1.045662289672538
$ws.Range("E22").Value = 1.052989484828921
$ws.Range("F22").Value = 1.060454435531443
$ws.Range("I22").Value = 1.044574987175653
$ws.Range("J22").Value = 1.051515410977896
$ws.Range("K22").Value = 1.049334448467771
$ws.Range("L22").Value = 1.056634116496416
$ws.Range("M22").Value = 1.06407145519093
$ws.Range("N22").Value = 1.053008683471988
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.045212442102327
$ws.Range("D23").Value = 1.046012550677569
$ws.Range("E23").Value = 1.053523519011335
$ws.Range("F23").Value = 1.061027191931126
$ws.Range("I23").Value = 1.044732211420059
$ws.Range("J23").Value = 1.051846020737517
$ws.Range("K23").Value = 1.04961584691086
$ws.Range("L23").Value = 1.057099176224751
$ws.Range("M23").Value = 1.064575668434208
$ws.Range("N23").Value = 1.053339762735389
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.047027767393311
$ws.Range("D24").Value = 1.047389315354819
$ws.Range("E24").Value = 1.055624965104079
$ws.Range("F24").Value = 1.063281356156508
$ws.Range("I24").Value = 1.045346954682172
$ws.Range("J24").Value = 1.053144063338876
$ws.Range("K24").Value = 1.050720107243539
$ws.Range("L24").Value = 1.058927941003041
$ws.Range("M24").Value = 1.066558896593693
$ws.Range("N24").Value = 1.05463964870613
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.049129420237974
$ws.Range("D25").Value = 1.048982543795075
$ws.Range("E25").Value = 1.05806172335074
$ws.Range("F25").Value = 1.065895911905025
$ws.Range("I25").Value = 1.046051616834288
$ws.Range("J25").Value = 1.05464314868589
$ws.Range("K25").Value = 1.051994212357805
$ws.Range("L25").Value = 1.061045895517195
$ws.Range("M25").Value = 1.068856768071878
$ws.Range("N25").Value = 1.056140862926284
